# Automatic update of files.
# Update the "Förändrad" (Changed) date in column C for rows 2-13
# from 2023-11-13 (45243) to 2023-11-14 (45244).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 13; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45243) {
        $cell.Value = 45244
    }
}
